# MitsosBarton2006Ex323 __C_Stationary generator (alpha_non_zero)
# "expermits todos no convexos menos el 5to"
#
# Rewrites the numeric coefficients baked into the Leader/Follower
# expression strings (and their derived evaluation points) across all
# the worksheets of the workbook. The underlying cell contents are all
# plain text (even the ones that look like numbers), so every write
# that could be mistaken for a number is forced to Text format first
# and the style pointer is put back to "Normal" afterwards so we don't
# leave a lingering custom format behind.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($Cell, $Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

# --- Restricciones_del_lider --------------------------------------------
$wsLider = $wb.Worksheets.Item("Restricciones_del_lider")

Set-TextValue $wsLider.Range("A2") "2.09 - x"
Set-TextValue $wsLider.Range("B2") "-3.09"
Set-TextValue $wsLider.Range("D2") "0.86"

Set-TextValue $wsLider.Range("A3") "-2.09 + x"
Set-TextValue $wsLider.Range("B3") "1.0899999999999999"
Set-TextValue $wsLider.Range("D3") "0.62"

Set-TextValue $wsLider.Range("A4") "41.02289999999999 + x - y - 9(x^2)"
Set-TextValue $wsLider.Range("B4") "-40.02289999999999"
Set-TextValue $wsLider.Range("D4") "0.58"

# --- Restricciones_del_follower -----------------------------------------
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $wsFollower.Range("A2") "98.95959999999998 - 32.083999999999996y + (-0.5 + x)*(y^2)"
Set-TextValue $wsFollower.Range("B2") "-98.95959999999998"
Set-TextValue $wsFollower.Range("D2") "0.69"
Set-TextValue $wsFollower.Range("E2") "-0.7000000000000001"
Set-TextValue $wsFollower.Range("F2") "-0.6"

Set-TextValue $wsFollower.Range("A3") "-3.686 + 0.97y"
Set-TextValue $wsFollower.Range("B3") "2.686"
Set-TextValue $wsFollower.Range("D3") "0.65"
# E3 stays "0" (unchanged)
Set-TextValue $wsFollower.Range("F3") "8.4"

Set-TextValue $wsFollower.Range("A4") "-671.9633333333334 + 176.66666666666669y"
Set-TextValue $wsFollower.Range("B4") "670.3333333333334"
Set-TextValue $wsFollower.Range("D4") "0.32"
Set-TextValue $wsFollower.Range("E4") "0"
Set-TextValue $wsFollower.Range("F4") "5.300000000000001"

# --- Punto_modificado -----------------------------------------------------
$wsPunto = $wb.Worksheets.Item("Punto_modificado")

Set-TextValue $wsPunto.Range("A2") "2.09"
Set-TextValue $wsPunto.Range("B2") "3.8"

# --- Vector_bf --------------------------------------------------------------
$wsBf = $wb.Worksheets.Item("Vector_bf")

Set-TextValue $wsBf.Range("A2") "-44.36383333333334"

# --- Vector_BF ----------------------------------------------------------------
$wsBF = $wb.Worksheets.Item("Vector_BF")

Set-TextValue $wsBF.Range("A2") "30.430279999999996"
Set-TextValue $wsBF.Range("A3") "-13.485826000000001"

# --- Vector_Alpha (A2 is a genuine number, not text) -----------------------
$wsAlpha = $wb.Worksheets.Item("Vector_Alpha")
$wsAlpha.Range("A2").Value = 0.03
